$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A25').Value = 'notation-functionalty;view-notation-keys-page'
$ws.Range('B25').Value = 'failed'
$ws.Range('C25').Value = '2021-01-01 15_39_05'
$ws.Range('D25').Value = 'chrome'

$ws.Range('A26').Value = 'notation-functionalty;create-a-new-notation-key'
$ws.Range('B26').Value = 'failed'
$ws.Range('C26').Value = '2021-01-01 15_39_33'
$ws.Range('D26').Value = 'chrome'

$ws.Range('A27').Value = 'notation-functionalty;edit-notation-key'
$ws.Range('B27').Value = 'failed'
$ws.Range('C27').Value = '2021-01-01 15_40_01'
$ws.Range('D27').Value = 'chrome'

$ws.Range('A28').Value = 'notation-functionalty;delete-the-notation-key'
$ws.Range('B28').Value = 'failed'
$ws.Range('C28').Value = '2021-01-01 15_40_26'
$ws.Range('D28').Value = 'chrome'

$ws.Range('A29').Value = 'notation-functionalty;view-notation-keys-page'
$ws.Range('B29').Value = 'failed'
$ws.Range('C29').Value = '2021-01-01 15_49_42'
$ws.Range('D29').Value = 'chrome'

$ws.Range('A30').Value = 'notation-functionalty;create-a-new-notation-key'
$ws.Range('B30').Value = 'failed'
$ws.Range('C30').Value = '2021-01-01 15_50_11'
$ws.Range('D30').Value = 'chrome'

$ws.Range('A31').Value = 'notation-functionalty;edit-notation-key'
$ws.Range('B31').Value = 'failed'
$ws.Range('C31').Value = '2021-01-01 15_50_33'
$ws.Range('D31').Value = 'chrome'

$ws.Range('A32').Value = 'notation-functionalty;delete-the-notation-key'
$ws.Range('B32').Value = 'failed'
$ws.Range('C32').Value = '2021-01-01 15_51_03'
$ws.Range('D32').Value = 'chrome'

$ws.Range('A33').Value = 'notation-functionalty;view-notation-keys-page'
$ws.Range('B33').Value = 'failed'
$ws.Range('C33').Value = '2021-01-01 15_53_11'
$ws.Range('D33').Value = 'chrome'

$ws.Range('A34').Value = 'create-new-attestations;user--should-be-able-to-create-attestations-and-view--successfully-validating-the-message'
$ws.Range('B34').Value = 'failed'
$ws.Range('C34').Value = '2021-01-01 15_54_37'
$ws.Range('D34').Value = 'chrome'

$ws.Range('A35').Value = 'notation-functionalty;view-notation-keys-page'
$ws.Range('B35').Value = 'failed'
$ws.Range('C35').Value = '2021-01-02 09_37_51'
$ws.Range('D35').Value = 'chrome'

$ws.Range('A36').Value = 'notation-functionalty;create-a-new-notation-key'
$ws.Range('B36').Value = 'failed'
$ws.Range('C36').Value = '2021-01-02 09_38_22'
$ws.Range('D36').Value = 'chrome'

$ws.Range('A37').Value = 'notation-functionalty;edit-notation-key'
$ws.Range('B37').Value = 'failed'
$ws.Range('C37').Value = '2021-01-02 09_38_51'
$ws.Range('D37').Value = 'chrome'

$ws.Range('A38').Value = 'notation-functionalty;delete-the-notation-key'
$ws.Range('B38').Value = 'failed'
$ws.Range('C38').Value = '2021-01-02 09_39_16'
$ws.Range('D38').Value = 'chrome'
